$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change "Obrigatorio" column (E) from "N" to "S" for rows 2 through 8
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 5).Value = "S"
}
